# Add a new "2025-advent" game entry into the "games" sheet, row 16.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("games")

$ws.Range("A16").Value = "2025-advent"
$ws.Range("B16").Value = "explore"
$ws.Range("C16").Value = "Kispest Kertváros - 2025 Adventi Ablakok"
$ws.Range("D16").Value = "6+"
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = "47.434879043711284, 19.1640427"

$ws.Activate()
$ws.Range("I16").Select()
